$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2024-05-27","11:28:58","-","-","Screw K30 no lo detecta puesto","-","-","11:29:04","0:00:06"),
    @("2024-05-27","11:29:03","-","-","No lee QR","-","-","11:29:04","0:00:01"),
    @("2024-05-27","11:29:06","-","-","Etiquetadora","-","-","11:29:08","0:00:02"),
    @("2024-05-27","11:29:24","-","-","La cámara no detecta Top Cover","-","-","11:29:26","0:00:02"),
    @("2024-05-27","11:29:27","-","-","No lee QR","-","-","11:29:28","0:00:01"),
    @("2024-05-27","11:29:29","-","-","Power atascado en prensa, cuesta sacar","-","-","11:29:31","0:00:02"),
    @("2024-05-27","11:29:37","-","-","La cámara no detecta Busbar","-","-","11:29:41","0:00:04"),
    @("2024-05-27","11:29:40","-","-","La cámara no detecta Top Cover","-","-","11:29:41","0:00:01"),
    @("2024-05-27","11:40:49","-","Cámara no detecta foams","-","-","-","11:40:52","0:00:03"),
    @("2024-05-27","11:40:53","-","Cámara no detecta Power CP","-","-","-","11:40:55","0:00:02"),
    @("2024-05-27","11:41:47","-","-","Atasco tuerca","-","-","11:41:52","0:00:05"),
    @("2024-05-27","11:41:50","-","-","No lee QR","-","-","11:41:51","0:00:01"),
    @("2024-05-27","11:41:53","-","-","La cámara no detecta Busbar","-","-","11:41:55","0:00:02")
)

$startRow = 287
$endRow = $startRow + $data.Length - 1

# Format column A as Text first so that date-looking strings like
# "2024-05-27" are stored as literal text instead of being auto-converted
# to date serial numbers.
$dateRange = $ws.Range("A$startRow`:A$endRow")
$dateRange.NumberFormat = "@"

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $rowValues = $data[$i]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowValues[$c]
    }
}

# Clear the temporary formatting so the new cells end up with the default
# (unstyled) look, matching the rest of the sheet.
$dateRange.ClearFormats()
